$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 231 (high/close revised) ---
$ws.Range("D231").Value = 458.92
$ws.Range("F231").Value = 451.56

# --- Append row 232 ---
$ws.Range("A231").Copy()
$ws.Range("A232").PasteSpecial(-4122)  # xlPasteFormats - reuse the datetime cell style
$ws.Range("A232").Value = 45047.33333333334
$ws.Range("B232").Value = "FX_IDC:USDKZT"
$ws.Range("C232").Value = 444.76
$ws.Range("D232").Value = 452.98
$ws.Range("E232").Value = 440.61
$ws.Range("F232").Value = 446.11
$ws.Range("G232").Value = 0

# --- Append row 233 ---
$ws.Range("A231").Copy()
$ws.Range("A233").PasteSpecial(-4122)
$ws.Range("A233").Value = 45078.33333333334
$ws.Range("B233").Value = "FX_IDC:USDKZT"
$ws.Range("C233").Value = 446.11
$ws.Range("D233").Value = 454.98
$ws.Range("E233").Value = 443.36
$ws.Range("F233").Value = 450.31
$ws.Range("G233").Value = 0

# --- Append row 234 ---
$ws.Range("A231").Copy()
$ws.Range("A234").PasteSpecial(-4122)
$ws.Range("A234").Value = 45110.33333333334
$ws.Range("B234").Value = "FX_IDC:USDKZT"
$ws.Range("C234").Value = 450.31
$ws.Range("D234").Value = 450.31
$ws.Range("E234").Value = 442.21
$ws.Range("F234").Value = 445.06
$ws.Range("G234").Value = 0

$excel.CutCopyMode = $false
